$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.606.80'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '1.853.44'
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '263.75'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5273'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3240'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06799'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.98'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7839'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07762'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").Value = '1.899.19'
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.73'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.037'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.99'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007965'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '26.627.78'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.640'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.477'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.016'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.81'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.170'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -6.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.679'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.01'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '111.95'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.184'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08721'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.104'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04861'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7226'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +6.22%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.133'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.875'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.111'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.269'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01792'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.4871'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9023'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.39'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.973'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.683'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05887'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.031'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1238'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.09'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.8906'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.14'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.67%  '
